# "El programa ejecuta un contador del 0 al 99"
# Adds a BCD-to-7-segment lookup sheet ("bcd_2_7seg") plus a small helper
# sheet ("Hoja3") and tweaks the timer-calculator sheet (Hoja1) so the
# computed OCR1A value drives a longer delay (Fclock/prescaler change).

$wb = $excel.ActiveWorkbook
$hoja1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Hoja1 tweaks: change Fclock and add the POWER(2,16) helper cell.
# ---------------------------------------------------------------------
$hoja1.Range("B1").Value = 0.1
$hoja1.Range("D6").Formula = "=POWER(2,16)"

# ---------------------------------------------------------------------
# 2) New sheets, inserted right after Hoja1, in this tab order:
#    Hoja1, Hoja3, bcd_2_7seg
# ---------------------------------------------------------------------
$hoja3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $hoja1)
$hoja3.Name = "Hoja3"

$bcd = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $hoja3)
$bcd.Name = "bcd_2_7seg"

# ---------------------------------------------------------------------
# 3) bcd_2_7seg: the BCD -> 7-segment switch table.
# ---------------------------------------------------------------------

# "Switch" cells (manual 1/0 inputs feeding the lookup row).
$bcd.Range("C2").Value = 1
$bcd.Range("B3").Value = 1
$bcd.Range("D3").Value = 1
$bcd.Range("C4").Value = 0
$bcd.Range("B5").Value = 1
$bcd.Range("D5").Value = 1
$bcd.Range("C6").Value = 1

# Header row (8) - segment labels (entered right-to-left: a,b,c,...,dot).
$bcd.Range("N8").Value = "a"
$bcd.Range("M8").Value = "b"
$bcd.Range("L8").Value = "c"
$bcd.Range("K8").Value = "d"
$bcd.Range("J8").Value = "e"
$bcd.Range("I8").Value = "f"
$bcd.Range("H8").Value = "g"
$bcd.Range("G8").Value = "dot"
$bcd.Range("O8").Value = "Binario"

# Row 9 - formulas pulling the switch states, then the binary string.
$bcd.Range("G9").Formula = "=F6"
$bcd.Range("H9").Formula = "=C4"
$bcd.Range("I9").Formula = "=B3"
$bcd.Range("J9").Formula = "=B5"
$bcd.Range("K9").Formula = "=C6"
$bcd.Range("L9").Formula = "=D5"
$bcd.Range("M9").Formula = "=D3"
$bcd.Range("N9").Formula = "=C2"
$bcd.Range("O9").Formula = '=CONCAT("0b"&G9&H9&I9&J9&K9&L9&M9&N9)'

# ---------------------------------------------------------------------
# 4) Visual styling of the switch block (B2:F6) - light gray fill plus
#    a best-effort reproduction of the per-cell border layout.
# ---------------------------------------------------------------------
$bcd.Range("B2:F6").Interior.ThemeColor = 2

# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10
# xlThin=2, xlMedium=-4138
$bcd.Range("B2").Borders.Item(7).Weight = -4138
$bcd.Range("B2").Borders.Item(8).Weight = -4138

$bcd.Range("C2").Borders.Item(7).Weight = 2
$bcd.Range("C2").Borders.Item(10).Weight = 2
$bcd.Range("C2").Borders.Item(8).Weight = -4138
$bcd.Range("C2").Borders.Item(9).Weight = 2

$bcd.Range("D2").Borders.Item(8).Weight = -4138
$bcd.Range("E2").Borders.Item(8).Weight = -4138

$bcd.Range("F2").Borders.Item(10).Weight = -4138
$bcd.Range("F2").Borders.Item(8).Weight = -4138

$bcd.Range("B3").Borders.Item(7).Weight = -4138
$bcd.Range("B3").Borders.Item(10).Weight = 2
$bcd.Range("B3").Borders.Item(8).Weight = 2
$bcd.Range("B3").Borders.Item(9).Weight = 2

$bcd.Range("D3,C4,D5").Borders.LineStyle = 1

$bcd.Range("F3,F4,F5").Borders.Item(10).Weight = -4138

$bcd.Range("B4").Borders.Item(7).Weight = -4138

$bcd.Range("B5").Borders.Item(7).Weight = -4138
$bcd.Range("B5").Borders.Item(10).Weight = 2
$bcd.Range("B5").Borders.Item(8).Weight = 2
$bcd.Range("B5").Borders.Item(9).Weight = 2

$bcd.Range("B6").Borders.Item(7).Weight = -4138
$bcd.Range("B6").Borders.Item(9).Weight = -4138

$bcd.Range("C6").Borders.Item(7).Weight = 2
$bcd.Range("C6").Borders.Item(10).Weight = 2
$bcd.Range("C6").Borders.Item(8).Weight = 2
$bcd.Range("C6").Borders.Item(9).Weight = -4138

$bcd.Range("D6").Borders.Item(9).Weight = -4138
$bcd.Range("E6").Borders.Item(9).Weight = -4138

$bcd.Range("F6").Borders.Item(7).Weight = 2
$bcd.Range("F6").Borders.Item(10).Weight = -4138
$bcd.Range("F6").Borders.Item(8).Weight = 2
$bcd.Range("F6").Borders.Item(9).Weight = -4138

# Column widths / row heights so the table reads like a little IC/DIP
# switch diagram.
$bcd.Columns.Item(1).ColumnWidth = 4.42578125
$bcd.Columns.Item(2).ColumnWidth = 3.5703125
$bcd.Columns.Item(3).ColumnWidth = 7.7109375
$bcd.Columns.Item(4).ColumnWidth = 4
$bcd.Columns.Item(5).ColumnWidth = 1.28515625
$bcd.Columns.Item(6).ColumnWidth = 2.85546875
$bcd.Range("G1:N1").ColumnWidth = 6.5703125
$bcd.Columns.Item(15).ColumnWidth = 11.85546875

$bcd.Rows.Item(1).RowHeight = 15.75
$bcd.Rows.Item(2).RowHeight = 16.5
$bcd.Rows.Item(3).RowHeight = 36.75
$bcd.Rows.Item(4).RowHeight = 16.5
$bcd.Rows.Item(5).RowHeight = 36.75
$bcd.Rows.Item(6).RowHeight = 13.5
$bcd.Rows.Item(7).RowHeight = 36.75
$bcd.Rows.Item(8).RowHeight = 18
$bcd.Rows.Item(9).RowHeight = 18

$bcd.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 5) Conditional formatting: switches show red when ON (1) and gray
#    when OFF (0).
# ---------------------------------------------------------------------
foreach ($addr in @("C2","D3","D5","C6","B5","C4","B3","F6")) {
    $cell = $bcd.Range($addr)

    $fcOn = $cell.FormatConditions.Add(1, 3, 1)
    $fcOn.Font.Color = 192
    $fcOn.Interior.Color = 192

    $fcOff = $cell.FormatConditions.Add(1, 3, 0)
    $fcOff.Font.Color = 8421504
    $fcOff.Interior.Color = 8421504
}

# --- Hoja3: a single reference note (typed last) ----------------------
$hoja3.Range("A1").Value = "admux register"

# ---------------------------------------------------------------------
# 6) Selections / active sheet.
# ---------------------------------------------------------------------
$hoja1.Range("B6").Select() | Out-Null
$bcd.Range("O9").Select() | Out-Null
$hoja3.Activate() | Out-Null
$hoja3.Range("A2").Select() | Out-Null
